$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26..106 down to 27..107.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with a new data record (same as the former
# row 26 record but with an updated date).
$ws.Cells.Item(26, 1).Value = 9
$ws.Cells.Item(26, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(26, 3).Value = "Metropolitana"
$ws.Cells.Item(26, 4).Value = 45099
$ws.Cells.Item(26, 5).Value = 13
$ws.Cells.Item(26, 6).Value = 100112029
$ws.Cells.Item(26, 7).Value = "Orégano"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 16
$ws.Cells.Item(26, 11).Value = 20000
$ws.Cells.Item(26, 12).Value = 20000
$ws.Cells.Item(26, 13).Value = 20000
$ws.Cells.Item(26, 14).Value = "`$/docena de atados"
$ws.Cells.Item(26, 15).Value = "Región Metropolitana"
$ws.Cells.Item(26, 16).Value = 6667
$ws.Cells.Item(26, 17).Value = 3
$ws.Cells.Item(26, 18).Value = "Hortaliza"
